# chore: update Sheets via scheduled runner
# Refreshes market-board derived columns (currentAveragePrice*, LevePrice*,
# LeveProfit*) on several leve rows across sheets, as pulled by the
# scheduled price-update job.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1465.1522
$ws.Range("I62").Value = 1437.6562
$ws.Range("J62").Value = 1528
$ws.Range("K62").Value = 1437.6562
$ws.Range("L62").Value = 1528
$ws.Range("M62").Value = -813.6561999999999
$ws.Range("N62").Value = -2776

$ws.Range("H65").Value = 1465.1522
$ws.Range("I65").Value = 1437.6562
$ws.Range("J65").Value = 1528
$ws.Range("K65").Value = 7188.280999999999
$ws.Range("L65").Value = 7640
$ws.Range("M65").Value = -4068.280999999999
$ws.Range("N65").Value = -13880

$ws.Range("H92").Value = 640.5333000000001
$ws.Range("I92").Value = 640.8
$ws.Range("J92").Value = 640
$ws.Range("K92").Value = 640.8
$ws.Range("L92").Value = 640
$ws.Range("M92").Value = 607.2
$ws.Range("N92").Value = -3136

$ws.Range("H96").Value = 826.5
$ws.Range("I96").Value = 435.33334
$ws.Range("J96").Value = 2000
$ws.Range("K96").Value = 1306.00002
$ws.Range("L96").Value = 6000
$ws.Range("M96").Value = 66.99998000000005
$ws.Range("N96").Value = -8746

$ws.Range("H131").Value = 2249.2856
$ws.Range("I131").Value = 598.75
$ws.Range("J131").Value = 4450
$ws.Range("K131").Value = 1796.25
$ws.Range("L131").Value = 13350
$ws.Range("M131").Value = 3243.75
$ws.Range("N131").Value = -23430

$ws.Range("H132").Value = 27303.205
$ws.Range("I132").Value = 30993.5
$ws.Range("J132").Value = 2209.2
$ws.Range("K132").Value = 92980.5
$ws.Range("L132").Value = 6627.599999999999
$ws.Range("M132").Value = -90450.5
$ws.Range("N132").Value = -11687.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1078.0541
$ws.Range("I97").Value = 1039.5
$ws.Range("K97").Value = 1039.5
$ws.Range("M97").Value = -543.5

$ws.Range("H110").Value = 986.6667
$ws.Range("I110").Value = 789.44446
$ws.Range("J110").Value = 1282.5
$ws.Range("K110").Value = 789.44446
$ws.Range("L110").Value = 1282.5
$ws.Range("M110").Value = 1255.55554
$ws.Range("N110").Value = -5372.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 20000
$ws.Range("J132").Value = 20000
$ws.Range("L132").Value = 20000
$ws.Range("N132").Value = -30120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1519.2593
$ws.Range("I99").Value = 1442.2858
$ws.Range("J99").Value = 1661.0526
$ws.Range("K99").Value = 1442.2858
$ws.Range("L99").Value = 1661.0526
$ws.Range("M99").Value = 55.71419999999989
$ws.Range("N99").Value = -4657.0526

$ws.Range("H126").Value = 1519.2593
$ws.Range("I126").Value = 1442.2858
$ws.Range("J126").Value = 1661.0526
$ws.Range("K126").Value = 4326.857400000001
$ws.Range("L126").Value = 4983.1578
$ws.Range("M126").Value = -1856.857400000001
$ws.Range("N126").Value = -9923.157800000001

$ws.Range("H132").Value = 75492.78999999999
$ws.Range("I132").Value = 3128.5715
$ws.Range("J132").Value = 147857
$ws.Range("K132").Value = 9385.7145
$ws.Range("L132").Value = 443571
$ws.Range("M132").Value = -6855.7145
$ws.Range("N132").Value = -448631

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4022.5
$ws.Range("I3").Value = 2030
$ws.Range("J3").Value = 10000
$ws.Range("K3").Value = 6090
$ws.Range("L3").Value = 30000
$ws.Range("M3").Value = -5978
$ws.Range("N3").Value = -30224

$ws.Range("H117").Value = 1157.0769
$ws.Range("I117").Value = 976.3333
$ws.Range("J117").Value = 1211.3
$ws.Range("K117").Value = 2928.9999
$ws.Range("L117").Value = 3633.9
$ws.Range("M117").Value = 513.0001000000002
$ws.Range("N117").Value = -10517.9

$ws.Range("H129").Value = 15153225
$ws.Range("I129").Value = 643.3333
$ws.Range("J129").Value = 17545738
$ws.Range("K129").Value = 1929.9999
$ws.Range("L129").Value = 52637214
$ws.Range("M129").Value = 3070.0001
$ws.Range("N129").Value = -52647214

$ws.Range("H131").Value = 6330077.5
$ws.Range("J131").Value = 6494413
$ws.Range("L131").Value = 19483239
$ws.Range("N131").Value = -19493319

$ws.Range("H133").Value = 6136.7676
$ws.Range("I133").Value = 2612.3845
$ws.Range("J133").Value = 7664
$ws.Range("K133").Value = 7837.1535
$ws.Range("L133").Value = 22992
$ws.Range("M133").Value = -2777.1535
$ws.Range("N133").Value = -33112

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2160.9167
$ws.Range("I97").Value = 1055
$ws.Range("J97").Value = 2713.875
$ws.Range("K97").Value = 1055
$ws.Range("L97").Value = 2713.875
$ws.Range("M97").Value = -559
$ws.Range("N97").Value = -3705.875

$ws.Range("H107").Value = 507.69232
$ws.Range("I107").Value = 410.77777
$ws.Range("J107").Value = 725.75
$ws.Range("K107").Value = 410.77777
$ws.Range("L107").Value = 725.75
$ws.Range("M107").Value = 1509.22223
$ws.Range("N107").Value = -4565.75

$ws.Range("H132").Value = 33214.03
$ws.Range("I132").Value = 1306.5238
$ws.Range("J132").Value = 94128.37
$ws.Range("K132").Value = 3919.5714
$ws.Range("L132").Value = 282385.11
$ws.Range("M132").Value = -1389.5714
$ws.Range("N132").Value = -287445.11

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1603.4
$ws.Range("I7").Value = 1558.7778
$ws.Range("J7").Value = 2005
$ws.Range("K7").Value = 1558.7778
$ws.Range("L7").Value = 2005
$ws.Range("M7").Value = -1446.7778
$ws.Range("N7").Value = -2229

$ws.Range("H16").Value = 1107.5
$ws.Range("I16").Value = 990
$ws.Range("J16").Value = 1225
$ws.Range("K16").Value = 990
$ws.Range("L16").Value = 1225
$ws.Range("M16").Value = -820
$ws.Range("N16").Value = -1565

$ws.Range("H93").Value = 2775.25
$ws.Range("I93").Value = 3219.125
$ws.Range("K93").Value = 3219.125
$ws.Range("M93").Value = -1971.125

$ws.Range("H119").Value = 29915.238
$ws.Range("J119").Value = 29915.238
$ws.Range("L119").Value = 29915.238
$ws.Range("N119").Value = -39591.238

$ws.Range("H122").Value = 2092.6428
$ws.Range("I122").Value = 1884.8
$ws.Range("J122").Value = 2612.25
$ws.Range("K122").Value = 5654.4
$ws.Range("L122").Value = 7836.75
$ws.Range("M122").Value = -3204.4
$ws.Range("N122").Value = -12736.75

$ws.Range("H126").Value = 1603.4
$ws.Range("I126").Value = 1558.7778
$ws.Range("J126").Value = 2005
$ws.Range("K126").Value = 4676.3334
$ws.Range("L126").Value = 6015
$ws.Range("M126").Value = -2206.3334
$ws.Range("N126").Value = -10955
